$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data changes: nodeId (B) becomes numeric, sound_type (E) values lose their quotes
$ws.Range("B2").Value = 1648804145
$ws.Range("B3").Value = 1648804146
$ws.Range("B4").Value = 1648804146

$ws.Range("E2").Value = "gun"
$ws.Range("E3").Value = "car"
$ws.Range("E4").Value = "animal"

# Widen the probability column (F) a bit, as in the new filter UI
$ws.Columns.Item(6).ColumnWidth = 9.83

# Selection moved from G4 to B2
$ws.Range("B2").Select() | Out-Null
